$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 5613
$ws.Range("E2").Value = 39
$ws.Range("F2").Value = 39
$ws.Range("G2").Value = -156
$ws.Range("H2").Value = -153
$ws.Range("I2").Value = -136
$ws.Range("J2").Value = -16
$ws.Range("K2").Value = 4524
$ws.Range("L2").Value = 2940
$ws.Range("M2").Value = 1584
$ws.Range("N2").Value = 1378
$ws.Range("O2").Value = 207
$ws.Range("P2").Value = 73
$ws.Range("Q2").Value = 445
$ws.Range("R2").Value = 159
$ws.Range("S2").Value = -966
$ws.Range("T2").Value = 213
$ws.Range("U2").Value = 232
$ws.Range("V2").Value = 1904
$ws.Range("W2").Value = 0.7
$ws.Range("X2").Value = -2.72
$ws.Range("Y2").Value = -9.57
$ws.Range("Z2").Value = -3.01
$ws.Range("AA2").Value = 185.57
$ws.Range("AB2").Value = 1840.75
$ws.Range("AC2").Value = -933
$ws.Range("AD2").Value = -3.13
$ws.Range("AE2").Value = 9953
$ws.Range("AF2").Value = 0.29
$ws.Range("AG2").Value = 100
$ws.Range("AH2").Value = 3.42
$ws.Range("AI2").Value = -10.15
$ws.Range("AJ2").Value = 14625466

$ws.Range("D3").Value = 5742
$ws.Range("E3").Value = 53
$ws.Range("F3").Value = 53
$ws.Range("G3").Value = 45
$ws.Range("H3").Value = 11
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 6
$ws.Range("K3").Value = 4547
$ws.Range("L3").Value = 2969
$ws.Range("M3").Value = 1579
$ws.Range("N3").Value = 1370
$ws.Range("O3").Value = 208
$ws.Range("P3").Value = 73
$ws.Range("Q3").Value = 243
$ws.Range("R3").Value = 37
$ws.Range("S3").Value = -163
$ws.Range("T3").Value = 256
$ws.Range("U3").Value = -13
$ws.Range("V3").Value = 1814
$ws.Range("W3").Value = 0.93
$ws.Range("X3").Value = 0.2
$ws.Range("Y3").Value = 0.39
$ws.Range("Z3").Value = 0.25
$ws.Range("AA3").Value = 188.05
$ws.Range("AB3").Value = 1822.1
$ws.Range("AC3").Value = 36
$ws.Range("AD3").Value = 96.74
$ws.Range("AE3").Value = 9900
$ws.Range("AF3").Value = 0.36
$ws.Range("AG3").Value = 100
$ws.Range("AH3").Value = 2.84
$ws.Range("AI3").Value = 260.11
$ws.Range("AJ3").Value = 14625466

$ws.Range("D4").Value = 5909
$ws.Range("E4").Value = 161
$ws.Range("F4").Value = 161
$ws.Range("G4").Value = 97
$ws.Range("H4").Value = 69
$ws.Range("I4").Value = 60
$ws.Range("J4").Value = 9
$ws.Range("K4").Value = 4552
$ws.Range("L4").Value = 2975
$ws.Range("M4").Value = 1577
$ws.Range("N4").Value = 1376
$ws.Range("O4").Value = 201
$ws.Range("P4").Value = 73
$ws.Range("Q4").Value = 334
$ws.Range("R4").Value = -339
$ws.Range("S4").Value = 14
$ws.Range("T4").Value = 330
$ws.Range("U4").Value = 4
$ws.Range("V4").Value = 1836
$ws.Range("W4").Value = 2.73
$ws.Range("X4").Value = 1.16
$ws.Range("Y4").Value = 4.37
$ws.Range("Z4").Value = 1.51
$ws.Range("AA4").Value = 188.69
$ws.Range("AB4").Value = 1875.38
$ws.Range("AC4").Value = 411
$ws.Range("AD4").Value = 8.83
$ws.Range("AE4").Value = 9938
$ws.Range("AF4").Value = 0.36
$ws.Range("AG4").Value = 150
$ws.Range("AH4").Value = 4.14
$ws.Range("AI4").Value = 34.58
$ws.Range("AJ4").Value = 14625466

$ws.Range("D5").Value = 5501
$ws.Range("E5").Value = -60
$ws.Range("F5").Value = -60
$ws.Range("G5").Value = -36
$ws.Range("H5").Value = -55
$ws.Range("I5").Value = -20
$ws.Range("J5").Value = -35
$ws.Range("K5").Value = 5090
$ws.Range("L5").Value = 3223
$ws.Range("M5").Value = 1867
$ws.Range("N5").Value = 1650
$ws.Range("O5").Value = 217
$ws.Range("P5").Value = 73
$ws.Range("Q5").Value = 165
$ws.Range("R5").Value = -223
$ws.Range("S5").Value = -3
$ws.Range("T5").Value = 345
$ws.Range("U5").Value = -181
$ws.Range("V5").Value = 1834
$ws.Range("W5").Value = -1.09
$ws.Range("X5").Value = -1
$ws.Range("Y5").Value = -1.32
$ws.Range("Z5").Value = -1.14
$ws.Range("AA5").Value = 172.59
$ws.Range("AB5").Value = 1828.19
$ws.Range("AC5").Value = -137
$ws.Range("AD5").Value = -25.08
$ws.Range("AE5").Value = 11920
$ws.Range("AF5").Value = 0.29
$ws.Range("AG5").Value = 150
$ws.Range("AH5").Value = 4.37
$ws.Range("AI5").Value = -103.68
$ws.Range("AJ5").Value = 14625466

$ws.Range("D6").Value = 5560
$ws.Range("E6").Value = -248
$ws.Range("F6").Value = -248
$ws.Range("G6").Value = -626
$ws.Range("H6").Value = -636
$ws.Range("I6").Value = -396
$ws.Range("K6").Value = 4954
$ws.Range("L6").Value = 3562
$ws.Range("M6").Value = 1392
$ws.Range("N6").Value = 1317
$ws.Range("P6").Value = 73
$ws.Range("Q6").Value = -217
$ws.Range("R6").Value = -190
$ws.Range("S6").Value = 370
$ws.Range("T6").Value = 259
$ws.Range("U6").Value = -476
$ws.Range("V6").Value = 2012
$ws.Range("W6").Value = -4.47
$ws.Range("X6").Value = -11.44
$ws.Range("Y6").Value = -26.68
$ws.Range("Z6").Value = -12.66
$ws.Range("AA6").Value = 255.96
$ws.Range("AB6").Value = 1192.98
$ws.Range("AC6").Value = -2706
$ws.Range("AD6").Value = -2.29
$ws.Range("AE6").Value = 9511
$ws.Range("AF6").Value = 0.65
$ws.Range("AG6").Value = 150
$ws.Range("AH6").Value = 2.42
$ws.Range("AI6").Value = -5.25
$ws.Range("AJ6").Value = 14625466

$ws.Range("D7:AJ9").ClearContents()

